# Apply the "adding more space in semillerostable" edit.
# Updates the activities sheet header info, the month columns (shift from
# Feb-Jun to Aug-Dec), clears one "X" mark, flips an activity status, and
# updates the "Integrantes" roster with new placeholder test data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1. Actividades")
$ws2 = $wb.Worksheets.Item("Hoja2. Integrantes")

# --- Sheet 1: header block ---
$ws1.Range("E4").Value = "Semillero de Investigación en Inteligencia Computacional"
$ws1.Range("E5").Value = "José Alejandro Cortés Taborda"
$ws1.Range("E6").Value = "2019-2"

# --- Sheet 1: month header row (row 9), shifted from Feb-Jun to Aug-Dec ---
$ws1.Range("D9").Value = "A"
$ws1.Range("E9").Value = "S"
$ws1.Range("F9").Value = "O"
$ws1.Range("G9").Value = "N"
$ws1.Range("H9").Value = "D"

# --- Sheet 1: activity row (row 10) ---
$ws1.Range("F10").ClearContents()
$ws1.Range("J10").Value = "Se realizó"

# --- Sheet 2: integrantes roster row (row 6) ---
$ws2.Range("A6").Value = "Usuario nuevo Prueba"
$ws2.Range("C6").Value = 1002
$ws2.Range("D6").Value = 1002
$ws2.Range("E6").Value = "Usuario_prueba@elpoli.edu.co"
